$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.438.04'
$ws.Range("E2").Value = '  +3.90%  '
$ws.Range("D3").Value = '2.446.78'
$ws.Range("E3").Value = '  +3.71%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = "'557.30"
$ws.Range("E5").Value = '  +2.86%  '
$ws.Range("D6").Value = "'139.21"
$ws.Range("E6").Value = '  +2.06%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = "'0.573"
$ws.Range("E8").Value = '  +1.55%  '
$ws.Range("E9").Value = '  +3.97%  '
$ws.Range("E10").Value = '  +4.16%  '
$ws.Range("E11").Value = '  +1.87%  '
$ws.Range("E12").Value = '  -2.10%  '
$ws.Range("D13").Value = "'25.01"
$ws.Range("E13").Value = '  +4.46%  '
$ws.Range("D14").Value = '2.879.71'
$ws.Range("E14").Value = '  +3.77%  '
$ws.Range("D15").Value = '60.352.39'
$ws.Range("E15").Value = '  +3.83%  '
$ws.Range("D16").Value = "'0.0000140"
$ws.Range("E16").Value = '  +4.74%  '
$ws.Range("D17").Value = '2.442.09'
$ws.Range("E17").Value = '  +3.89%  '
$ws.Range("D18").Value = "'11.50"
$ws.Range("E18").Value = '  +7.03%  '
$ws.Range("E19").Value = '  +3.53%  '
$ws.Range("D20").Value = "'335.98"
$ws.Range("E20").Value = '  +0.98%  '
$ws.Range("E21").Value = '  +2.27%  '
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = "'64.71"
$ws.Range("E23").Value = '  +2.83%  '
$ws.Range("E24").Value = '  +2.17%  '
$ws.Range("D25").Value = "'8.56"
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = '  +0.34%  '
$ws.Range("E27").Value = '  -0.52%  '
$ws.Range("D28").Value = '0.0₃0798'
$ws.Range("E28").Value = '  +7.74%  '
$ws.Range("E29").Value = '  +3.41%  '
$ws.Range("E30").Value = '  +2.47%  '
$ws.Range("D31").Value = "'170.92"
$ws.Range("E31").Value = '  -0.99%  '
$ws.Range("E32").Value = '  +1.85%  '
$ws.Range("E33").Value = '  -1.01%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("E35").Value = '  +5.32%  '
$ws.Range("D36").Value = "'4.29"
$ws.Range("E36").Value = '  +1.07%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("D39").Value = "'40.09"
$ws.Range("E39").Value = '  +1.91%  '
$ws.Range("D40").Value = "'0.419"
$ws.Range("E40").Value = '  +10.78%  '
$ws.Range("D41").Value = "'317.09"
$ws.Range("E41").Value = '  +7.91%  '
$ws.Range("D42").Value = "'3.74"
$ws.Range("E42").Value = '  +2.38%  '
$ws.Range("D43").Value = "'144.30"
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("D44").Value = "'0.0965"
$ws.Range("E44").Value = '  +1.86%  '
$ws.Range("D45").Value = "'19.94"
$ws.Range("E45").Value = '  +2.96%  '
$ws.Range("E46").Value = '  +4.46%  '
$ws.Range("D47").Value = "'0.574"
$ws.Range("E47").Value = '  +1.43%  '
$ws.Range("D48").Value = "'0.406"
$ws.Range("E48").Value = '  +5.46%  '
$ws.Range("E49").Value = '  +2.64%  '
$ws.Range("E50").Value = '  -0.23%  '
$ws.Range("D51").Value = "'1.66"
$ws.Range("E51").Value = '  +6.03%  '
